$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "1.00", "5.74") are stored as text, not auto-converted to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.936.27"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.670.57"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "215.43"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.0621"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "20.31"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "1.905.84"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "1.683.43"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "65.63"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "26.937.54"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "234.69"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "7.95"
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "9.15"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").Value = "146.39"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "7.11"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "15.88"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.111"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "1.444.44"
$ws.Range("E33").Value = "  -5.40%  "
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "0.904"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +14.73%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.74"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.30"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "66.43"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "1.813.21"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").Value = "90.59"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +4.99%  "
$ws.Range("D50").Value = "0.0507"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  +0.22%  "
